$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 2.4
$ws.Range("R2").Value = 1.53
